$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '87.193.02'
$ws.Range("E2").Value = '  -3.07%  '
$ws.Range("D3").Value = '3.022.41'
$ws.Range("E3").Value = '  -6.49%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '205.04'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -6.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '608.68'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.44%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.356'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -9.47%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.800'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +14.19%  '
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("D10").Value = '3.021.98'
$ws.Range("E10").Value = '  -6.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.579'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.05%  '
$ws.Range("E12").Value = '  -1.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000228'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -12.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.16'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -4.35%  '
$ws.Range("D15").Value = '87.239.91'
$ws.Range("E15").Value = '  -2.54%  '
$ws.Range("D16").Value = '3.589.36'
$ws.Range("E16").Value = '  -5.81%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.74'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -7.92%  '
$ws.Range("D18").Value = '3.047.84'
$ws.Range("E18").Value = '  -4.82%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.11'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -7.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000194'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -18.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.77'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -5.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '411.57'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -5.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.90'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -8.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.73'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -7.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.22'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.92%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '79.61'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.19'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -5.50%  '
$ws.Range("D28").Value = '3.227.18'
$ws.Range("E28").Value = '  -4.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.07'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +7.27%  '
$ws.Range("E31").Value = '  -2.74%  '
$ws.Range("E32").Value = '  -7.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '492.26'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -9.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.39'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -17.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.76'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -8.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.40'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -8.84%  '
$ws.Range("E37").Value = '  -8.07%  '
$ws.Range("B38").Value = 'WhiteBITCoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.15'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.99%  '
$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '21.69'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.57%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.129'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.45%  '
$ws.Range("E41").Value = '  +0.37%  '
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '148.12'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.22%  '
$ws.Range("B44").Value = 'PolygonEcosystemToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.351'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -6.57%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.132'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +6.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.75'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -9.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.15'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.30%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0650'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +7.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '152.05'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -12.39%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.683'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -9.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.14'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -9.31%  '
